$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Worksheet, $CellRef, $Text)
    $cell = $Worksheet.Range($CellRef)
    $origStyle = $cell.Style
    # Leading apostrophe forces Excel to treat the value as literal text,
    # so numeric-looking strings (e.g. "1.00", "61.255.31") are not coerced
    # into numbers. Restore the original style afterward since the quote-prefix
    # write can tag the cell with a new style index.
    $cell.Value = "'" + $Text
    $cell.Style = $origStyle
}

Set-TextCell $ws "D2" "61.243.13"
Set-TextCell $ws "E2" "  +7.55%  "
Set-TextCell $ws "D3" "3.341.49"
Set-TextCell $ws "E3" "  +2.37%  "
Set-TextCell $ws "E4" "  +0.08%  "
Set-TextCell $ws "D5" "411.21"
Set-TextCell $ws "E5" "  +3.80%  "
Set-TextCell $ws "D6" "115.93"
Set-TextCell $ws "E6" "  +6.49%  "
Set-TextCell $ws "D7" "3.334.51"
Set-TextCell $ws "E7" "  +2.18%  "
Set-TextCell $ws "D8" "0.574"
Set-TextCell $ws "E8" "  -2.17%  "
Set-TextCell $ws "D9" "1.00"
Set-TextCell $ws "E9" "  +0.05%  "
Set-TextCell $ws "D10" "0.629"
Set-TextCell $ws "E10" "  +0.16%  "
Set-TextCell $ws "E11" "  +18.00%  "
Set-TextCell $ws "D12" "40.10"
Set-TextCell $ws "E12" "  +1.95%  "
Set-TextCell $ws "E13" "  -0.71%  "
Set-TextCell $ws "D14" "3.870.05"
Set-TextCell $ws "E14" "  +2.44%  "
Set-TextCell $ws "D15" "8.30"
Set-TextCell $ws "E15" "  -1.17%  "
Set-TextCell $ws "D16" "19.21"
Set-TextCell $ws "E16" "  +0.13%  "
Set-TextCell $ws "D17" "3.338.29"
Set-TextCell $ws "E17" "  +2.17%  "
Set-TextCell $ws "D18" "61.102.00"
Set-TextCell $ws "E18" "  +7.45%  "
Set-TextCell $ws "E19" "  -2.61%  "
Set-TextCell $ws "D20" "10.84"
Set-TextCell $ws "E20" "  +0.53%  "
Set-TextCell $ws "E21" "  +5.88%  "
Set-TextCell $ws "D22" "3.36"
Set-TextCell $ws "E22" "  +0.29%  "
Set-TextCell $ws "D23" "12.51"
Set-TextCell $ws "E23" "  -4.13%  "
Set-TextCell $ws "D24" "295.54"
Set-TextCell $ws "E24" "  +0.00%  "
Set-TextCell $ws "D25" "74.20"
Set-TextCell $ws "E25" "  -0.13%  "
Set-TextCell $ws "D26" "3.12"
Set-TextCell $ws "E26" "  -1.81%  "
Set-TextCell $ws "D27" "29.09"
Set-TextCell $ws "E27" "  +3.10%  "
Set-TextCell $ws "D28" "7.79"
Set-TextCell $ws "E28" "  +7.10%  "
Set-TextCell $ws "E29" "  -2.68%  "
Set-TextCell $ws "E30" "  +1.76%  "
Set-TextCell $ws "D31" "7.54"
Set-TextCell $ws "E31" "  -1.75%  "
Set-TextCell $ws "D32" "0.114"
Set-TextCell $ws "E32" "  +5.08%  "
Set-TextCell $ws "D33" "42.89"
Set-TextCell $ws "E33" "  +6.80%  "
Set-TextCell $ws "B34" "Dai"
Set-TextCell $ws "C34" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws "D34" "1.00"
Set-TextCell $ws "E34" "  +0.03%  "
Set-TextCell $ws "B35" "Toncoin"
Set-TextCell $ws "C35" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws "D35" "2.53"
Set-TextCell $ws "E35" "  +18.73%  "
Set-TextCell $ws "B36" "Cosmos"
Set-TextCell $ws "C36" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws "D36" "11.30"
Set-TextCell $ws "E36" "  +0.23%  "
Set-TextCell $ws "D37" "0.0489"
Set-TextCell $ws "E37" "  -0.97%  "
Set-TextCell $ws "D38" "52.46"
Set-TextCell $ws "E38" "  +1.82%  "
Set-TextCell $ws "E39" "  -0.09%  "
Set-TextCell $ws "D40" "3.06"
Set-TextCell $ws "E40" "  +4.35%  "
Set-TextCell $ws "D41" "3.43"
Set-TextCell $ws "E41" "  -1.72%  "
Set-TextCell $ws "D42" "135.09"
Set-TextCell $ws "E42" "  -2.92%  "
Set-TextCell $ws "E43" "  -1.31%  "
Set-TextCell $ws "D44" "0.288"
Set-TextCell $ws "E44" "  +1.94%  "
Set-TextCell $ws "D45" "1.90"
Set-TextCell $ws "E45" "  -0.30%  "
Set-TextCell $ws "D46" "3.88"
Set-TextCell $ws "E46" "  -3.34%  "
Set-TextCell $ws "D47" "16.36"
Set-TextCell $ws "E47" "  -4.59%  "
Set-TextCell $ws "E48" "  +4.36%  "
Set-TextCell $ws "D49" "21.19"
Set-TextCell $ws "E49" "  -4.83%  "
Set-TextCell $ws "D50" "2.151.98"
Set-TextCell $ws "E50" "  -0.69%  "
Set-TextCell $ws "D51" "3.671.05"
Set-TextCell $ws "E51" "  +2.55%  "
